# fix parameterization, run strategy and influence experiments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Update the "no estimate / default" row (row 30, cols C:G) from -0.5 to -0.1 ---
$ws.Range("C30:G30").Value = -0.1

# --- Update the view: scroll the frozen (bottom-right) pane so the window shows
#     further down the sheet, and leave the final selection at H36 ---
# Keep the freeze at row 1 / col 1 (header row + label column) exactly as before,
# just move the visible window/selection.
$ws.Range("B23").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 2

$ws.Range("H36").Select() | Out-Null
